$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Assign "Brian Galok" to D2, D4, D7 (Assigned to column)
$ws.Range("D2").Value = "Brian Galok"
$ws.Range("D4").Value = "Brian Galok"
$ws.Range("D7").Value = "Brian Galok"

# Reset the view: scroll back to top-left / unfreeze scroll position, select E2
$ws.Range("E2").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
